$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 189, pushing the existing rows 189:316 down to 190:317.
$ws.Rows(189).Insert()

# Populate the new row 189 with the new weekly data point.
$ws.Cells.Item(189, 1).Value = 10
$ws.Cells.Item(189, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(189, 3).Value = "La Araucanía"
$ws.Cells.Item(189, 4).Value = 45126
$ws.Cells.Item(189, 5).Value = 9
$ws.Cells.Item(189, 6).Value = 100112005
$ws.Cells.Item(189, 7).Value = "Puerro"
$ws.Cells.Item(189, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 50
$ws.Cells.Item(189, 11).Value = 8000
$ws.Cells.Item(189, 12).Value = 8000
$ws.Cells.Item(189, 13).Value = 8000
$ws.Cells.Item(189, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(189, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(189, 16).Value = 667
$ws.Cells.Item(189, 17).Value = 12
$ws.Cells.Item(189, 18).Value = "Hortaliza"
